$wb = $excel.ActiveWorkbook

# --- "adults" sheet (sheet2): append two rows copied from "peds" (sheet1) ---
$peds = $wb.Worksheets.Item("peds")
$adults = $wb.Worksheets.Item("adults")

# Row 11 <- peds row 2 (Hepatitis B)
$adults.Range("A11:O11").Value2 = $peds.Range("A2:O2").Value2

# Row 12 <- peds row 3 (Rotavirus)
$adults.Range("A12:M12").Value2 = $peds.Range("A3:M3").Value2

# --- selections as left by the author when the file was saved ---
$peds.Range("A2:S3").Select()

$adults.Select()
$adults.Range("B25").Select()
